# Updates the cryptos price table with the latest scraped prices/changes.
# Cells whose new value looks like a plain number (e.g. "336.16") are
# written with a leading apostrophe so Excel stores them as text instead
# of silently converting them to a numeric value (which would drop
# significant trailing/leading zeros and change the cell type). The
# Style is then reset to "Normal" so the transient quote-prefix style
# Excel applies doesn't change the cell's style index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.143.82'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.801.08'
$ws.Range("E3").Value = '  +2.86%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '''336.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").Value = '''0.4613'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +22.09%  '
$ws.Range("D8").Value = '''0.3710'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.58%  '
$ws.Range("D9").Value = '''45.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").Value = '''1.155'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.83%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '''0.07644'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.00%  '
$ws.Range("D12").Value = '''22.54'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = '''6.375'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.64%  '
$ws.Range("D15").Value = '''7.401'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.89%  '
$ws.Range("D16").Value = '1.799.79'
$ws.Range("E16").Value = '  +2.49%  '
$ws.Range("D17").Value = '''0.00001095'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.70%  '
$ws.Range("D18").Value = '''0.06737'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.29%  '
$ws.Range("D19").Value = '''83.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.50%  '
$ws.Range("D20").Value = '''1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '''17.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.44%  '
$ws.Range("D22").Value = '''6.433'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.96%  '
$ws.Range("D23").Value = '28.139.23'
$ws.Range("E23").Value = '  +0.32%  '
$ws.Range("E24").Value = '  +2.11%  '
$ws.Range("D25").Value = '''2.407'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("D26").Value = '''20.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.70%  '
$ws.Range("D27").Value = '''2.394'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.91%  '
$ws.Range("D28").Value = '''152.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("D29").Value = '2.004.53'
$ws.Range("E29").Value = '  +2.45%  '
$ws.Range("D30").Value = '''134.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.80%  '
$ws.Range("D31").Value = '''1.268'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("D32").Value = '''4.054'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("E33").Value = '  +11.21%  '
$ws.Range("D34").Value = '''5.899'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.98%  '
$ws.Range("B35").Value = 'Algorand'
$ws.Range("C35").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D35").Value = '''0.2240'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.21%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '''0.02381'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.80%  '
$ws.Range("D37").Value = '''12.22'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.06393'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.39%  '
$ws.Range("D39").Value = '''0.6726'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = '''5.280'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.52%  '
$ws.Range("D41").Value = '''1.528'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.53%  '
$ws.Range("D42").Value = '''1.233'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("D43").Value = '''8.126'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.74%  '
$ws.Range("D44").Value = '''14.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.58%  '
$ws.Range("D45").Value = '''1.000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '''0.6173'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.13%  '
$ws.Range("D47").Value = '''3.844'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.42%  '
$ws.Range("D48").Value = '''130.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("D49").Value = '''2.067'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.53%  '
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("D51").Value = '''0.07135'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.22%  '
